$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells we touch stay as text, matching original inlineStr type
$priceCells = @("D2","D3","D5","D6","D10","D11","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D24","D25","D26","D27","D29","D30","D31","D33","D34","D36","D38","D40","D42","D45","D46","D49","D50","D51","D8","D9")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "76.359.57"
$ws.Range("E2").Value = "  -0.25%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.075.86"
$ws.Range("E3").Value = "  +4.72%  "

# Row 5 - Solana
$ws.Range("D5").Value = "198.26"
$ws.Range("E5").Value = "  -0.29%  "

# Row 6 - BNB
$ws.Range("D6").Value = "618.38"
$ws.Range("E6").Value = "  +3.95%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.05%  "

# Row 8 - was Dogecoin, becomes XRP
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").Value = "0.552"
$ws.Range("E8").Value = "  +0.67%  "

# Row 9 - was XRP, becomes Dogecoin
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "0.214"
$ws.Range("E9").Value = "  +7.84%  "

# Row 10 - LidoStakedEther
$ws.Range("D10").Value = "3.075.91"
$ws.Range("E10").Value = "  +5.03%  "

# Row 11 - Cardano
$ws.Range("D11").Value = "0.445"
$ws.Range("E11").Value = "  +0.86%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +0.03%  "

# Row 13 - Toncoin
$ws.Range("D13").Value = "5.25"
$ws.Range("E13").Value = "  +7.29%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "3.638.45"
$ws.Range("E14").Value = "  +4.90%  "

# Row 15 - Avalanche
$ws.Range("D15").Value = "29.30"
$ws.Range("E15").Value = "  +2.93%  "

# Row 16 - ShibaInu
$ws.Range("D16").Value = "0.0000196"
$ws.Range("E16").Value = "  +3.48%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "76.233.89"
$ws.Range("E17").Value = "  -0.30%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.079.24"
$ws.Range("E18").Value = "  +5.45%  "

# Row 19 - Chainlink
$ws.Range("D19").Value = "13.50"
$ws.Range("E19").Value = "  -0.09%  "

# Row 20 - Uniswap
$ws.Range("D20").Value = "9.02"
$ws.Range("E20").Value = "  +3.12%  "

# Row 21 - was SuiNetwork, becomes BitcoinCash
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "384.25"
$ws.Range("E21").Value = "  +2.63%  "

# Row 22 - was BitcoinCash, becomes SuiNetwork
$ws.Range("B22").Value = "SuiNetwork"
$ws.Range("C22").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D22").Value = "2.58"
$ws.Range("E22").Value = "  +13.93%  "

# Row 23 - Polkadot
$ws.Range("E23").Value = "  +4.68%  "

# Row 24 - LEO
$ws.Range("D24").Value = "6.45"
$ws.Range("E24").Value = "  +0.78%  "

# Row 25 - NEARProtocol
$ws.Range("D25").Value = "4.57"
$ws.Range("E25").Value = "  +6.93%  "

# Row 26 - WrappedeETH
$ws.Range("D26").Value = "3.236.25"
$ws.Range("E26").Value = "  +4.79%  "

# Row 27 - Litecoin
$ws.Range("D27").Value = "72.42"
$ws.Range("E27").Value = "  +0.78%  "

# Row 28 - Dai
$ws.Range("E28").Value = "  +0.13%  "

# Row 29 - Aptos
$ws.Range("D29").Value = "10.08"
$ws.Range("E29").Value = "  +4.42%  "

# Row 30 - PEPE
$ws.Range("D30").Value = "0.0000108"
$ws.Range("E30").Value = "  +0.42%  "

# Row 31 - Binance-PegBSC-USD
$ws.Range("D31").Value = "0.993"
$ws.Range("E31").Value = "  -0.60%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("E32").Value = "  -0.43%  "

# Row 33 - Fetch.AI
$ws.Range("D33").Value = "1.42"
$ws.Range("E33").Value = "  +3.57%  "

# Row 34 - Bittensor
$ws.Range("D34").Value = "500.37"
$ws.Range("E34").Value = "  +0.09%  "

# Row 35 - PancakeSwap
$ws.Range("E35").Value = "  +5.67%  "

# Row 36 - Kaspa
$ws.Range("D36").Value = "0.129"
$ws.Range("E36").Value = "  +16.16%  "

# Row 37 - FirstDigitalUSD
$ws.Range("E37").Value = "  +0.01%  "

# Row 38 - EthereumClassic
$ws.Range("D38").Value = "20.88"
$ws.Range("E38").Value = "  +3.56%  "

# Row 39 - Monero
$ws.Range("E39").Value = "  -1.21%  "

# Row 40 - Aave
$ws.Range("D40").Value = "194.38"
$ws.Range("E40").Value = "  +8.44%  "

# Row 41 - WhiteBITCoin
$ws.Range("E41").Value = "  +0.56%  "

# Row 42 - PolygonEcosystemToken
$ws.Range("D42").Value = "0.378"
$ws.Range("E42").Value = "  -3.63%  "

# Row 43 - Cronos
$ws.Range("E43").Value = "  -7.14%  "

# Row 44 - USDe
$ws.Range("E44").Value = "  +0.06%  "

# Row 45 - Mantle
$ws.Range("D45").Value = "0.799"
$ws.Range("E45").Value = "  +20.62%  "

# Row 46 - RenderToken
$ws.Range("D46").Value = "5.16"
$ws.Range("E46").Value = "  +4.89%  "

# Row 47 - ImmutableX
$ws.Range("E47").Value = "  +6.53%  "

# Row 48 - Stacks
$ws.Range("E48").Value = "  +1.00%  "

# Row 49 - dogwifhat
$ws.Range("D49").Value = "2.45"
$ws.Range("E49").Value = "  +5.22%  "

# Row 50 - OKB
$ws.Range("D50").Value = "40.82"
$ws.Range("E50").Value = "  +2.14%  "

# Row 51 - ARBITRUM
$ws.Range("D51").Value = "0.598"
$ws.Range("E51").Value = "  +0.94%  "
